$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 73 is brand new (the sheet previously ended at row 72), so give it
# the same font as the rest of the data rows before populating it.
$ws.Range("A73:C73").Font.Name = "TAHOMA"

# Ensure CodiceFiscale column keeps leading zeros / is treated as text.
# Pre-format the exact target range (including the new row 73) as Text
# *before* writing any values, so Excel does not auto-coerce the
# numeric-looking CodiceFiscale strings into numbers (which would drop
# leading zeros).
$ws.Range("B2:B73").NumberFormat = "@"

$data = @(
    @('Ministero dell''Istruzione, dell''Universita'' e della Ricerca', '80185250588', 8457.0),
    @('Unione Italiana delle Camere di Commercio Industria, Artigianato e Agricoltura', '01484460587', 1051.0),
    @('Poste Italiane S.P.A.', '97103880585', 739.0),
    @('E-Fil S.r.l.', '03789750100', 472.0),
    @('Regione del Veneto', '80007580279', 402.0),
    @('CREDEMTEL SpA', '01378570350', 400.0),
    @('Regione Lombardia', '80050050154', 375.0),
    @('Lepida Spa', '02770891204', 372.0),
    @('Maggioli SPA', '06188330150', 371.0),
    @('Provincia Autonoma di Trento', '00337460224', 331.0),
    @('Banca Popolare di Sondrio, Società Cooperativa per Azioni', '00053810149', 311.0),
    @('Dedagroup Public Services S.R.L.', '01727860221', 232.0),
    @('Alto Adige Riscossioni Spa', '02805390214', 180.0),
    @('Progetti e Soluzioni SPA', '06423240727', 166.0),
    @('P.A. Digitale spa', '06628860964', 146.0),
    @('APKAPPA S.R.L.', '08543640158', 141.0),
    @('Regione Marche', '80008630420', 139.0),
    @('ADVANCED SYSTEMS srl', '03383350638', 126.0),
    @('Regione Puglia', '80017210727', 117.0),
    @('Regione Basilicata', '80002950766', 106.0),
    @('Regione Toscana', '01386030488', 100.0),
    @('ANCITEL', '07196850585', 90.0),
    @('Regione Piemonte', '80087670016', 90.0),
    @('Next Step Solution', '02554480349', 74.0),
    @('Regione Autonoma Friuli-Venezia Giulia', '80014930327', 73.0),
    @('Intesa Sanpaolo SPA', '00799960158', 59.0),
    @('NORDCOM', '13384100155', 55.0),
    @('Regione Umbria', '80000130544', 51.0),
    @('Bluenext S.r.l.', '04228480408', 49.0),
    @('Regione Autonoma della Sardegna', '80002870923', 47.0),
    @('Italriscossioni Società Italiana di Fiscalità Locale S.r.l.', '06092371001', 44.0),
    @('CINECA consorzio universitario', '00317740371', 41.0),
    @('Consorzio I.T. Srl', '01321400192', 40.0),
    @('UNIMATICA S.P.A', '02098391200', 38.0),
    @('ROMA CAPITALE', '02438750586', 37.0),
    @('SI.net Servizi Informatici S.r.L.', '02743730125', 35.0),
    @('PMPay s.r.l.', '08747230962', 33.0),
    @('Unicredit, Societa'' per Azioni', '00348170101', 33.0),
    @('DCS SOFTWARE E SERVIZI S.R.L.', '08063140019', 25.0),
    @('Regione Liguria', '00849050109', 23.0),
    @('Siscom SPA', '01778000040', 21.0),
    @('Regione Lazio', '80143490581', 18.0),
    @('Nexi SpA', '13212880150', 18.0),
    @('Citta'' Metropolitana di Roma Capitale', '80034390585', 18.0),
    @('Comune di Palermo', '80016350821', 17.0),
    @('ANDREANI TRIBUTI srl', '01412920439', 17.0),
    @('Novares Spa', '12105121003', 15.0),
    @('Servizi Locali SpA', '03170580751', 14.0),
    @('Regione Autonoma Valle D''Aosta', '80002270074', 13.0),
    @('Crédit Agricole Group Solutions Società Consortile per azioni', '02771790348', 10.0),
    @('Aric Agenzia Regionale di Informatica e Committenza', '91022630676', 10.0),
    @('UBI Banca', '03053920165', 8.0),
    @('Numera Sistemi e Informatica SpA', '01265230902', 8.0),
    @('Si.Form Consulting srl', '03943960827', 8.0),
    @('Be Smart s.r.l.', '05817461006', 6.0),
    @('Comune di Catania', '00137020871', 6.0),
    @('ISWEB S.p.A.', '01722270665', 3.0),
    @('Linea Comune Spa', '05591950489', 3.0),
    @('ARCA Servizi s.r.l', '09106071005', 2.0),
    @('Banco BPM Società per Azioni', '09722490969', 1.0),
    @('ARGO SOFTWARE SRL', '00838520880', 1.0),
    @('Società Almaviva S.p.A.', '08450891000', 1.0),
    @('MegASP S.r.l.', '09898030151', 1.0),
    @('San Marco SPA', '04142440728', 1.0),
    @('I.C.A. - Imposte Comunali Affini – s.r.l.', '02478610583', 1.0),
    @('e-SED Società Cooperativa', '02695640421', 1.0),
    @('Ministero dello Sviluppo Economico', '80230390587', 1.0),
    @('Softline srl', '12299030150', 1.0),
    @('CityPoste Payment Digital S.r.l.', '02003750672', 1.0),
    @('ICCREA Banca SpA', '04774801007', 1.0),
    @('Engineering Ingegneria Informatica SpA', '00967720285', 1.0),
    @('Agenzia Italiana del Farmaco - AIFA', '97345810580', 1.0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
